$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire item row for "مناديل شمع" (row 46). This shifts the
# totals row (47 -> 46) and the footer row (48 -> 47) up by one.
$ws.Rows.Item(46).Delete()

# Update the running total (column P) to reflect removal of the 40.00 item.
$ws.Range("P46").Value = 2449.65

# Update the generated timestamp shown in the footer.
$ws.Range("A47").Value = "Thursday, 25 September, 2025 4:53 PM"
